# Update the term-list example sheet to the latest Figure6 / OBCS release:
#  - the glossary entry for OBCS_0000120 is renamed from "set of data set"
#    to "data matrix"
#  - leave the last-used selection on B8 (matches the refreshed file)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "data matrix"

$ws.Range("B8").Select() | Out-Null
